{"js": "// Replace the date line and every \"a\u00d7b=c\" answer cell with the new values,\n// per the commit's diff. Every old string in the mapping is unique within\n// the document, so an exact, case-sensitive, non-wildcard search is safe.\nconst replacements = [\n  [\"2025-06-09 Monday\", \"2025-06-10 Tuesday\"],\n  [\"814\u00d79=7326\", \"880\u00d73=2640\"],\n  [\"259\u00d76=1554\", \"831\u00d78=6648\"],\n  [\"407\u00d78=3256\", \"687\u00d72=1374\"],\n  [\"541\u00d75=2705\", \"621\u00d76=3726\"],\n  [\"974\u00d79=8766\", \"474\u00d79=4266\"],\n  [\"390\u00d77=2730\", \"834\u00d74=3336\"],\n  [\"510\u00d76=3060\", \"707\u00d72=1414\"],\n  [\"736\u00d72=1472\", \"414\u00d78=3312\"],\n  [\"441\u00d77=3087\", \"861\u00d76=5166\"],\n  [\"205\u00d78=1640\", \"950\u00d75=4750\"],\n  [\"919\u00d73=2757\", \"408\u00d75=2040\"],\n  [\"843\u00d76=5058\", \"628\u00d73=1884\"],\n  [\"334\u00d79=3006\", \"801\u00d72=1602\"],\n  [\"946\u00d77=6622\", \"755\u00d77=5285\"],\n  [\"763\u00d76=4578\", \"575\u00d73=1725\"],\n  [\"864\u00d75=4320\", \"797\u00d75=3985\"],\n  [\"674\u00d75=3370\", \"533\u00d76=3198\"],\n  [\"388\u00d73=1164\", \"750\u00d78=6000\"],\n  [\"326\u00d72=652\", \"394\u00d73=1182\"],\n  [\"498\u00d75=2490\", \"851\u00d78=6808\"],\n  [\"125\u00d72=250\", \"317\u00d74=1268\"],\n  [\"201\u00d72=402\", \"983\u00d79=8847\"],\n  [\"231\u00d76=1386\", \"218\u00d76=1308\"],\n  [\"252\u00d77=1764\", \"261\u00d73=783\"],\n  [\"856\u00d75=4280\", \"300\u00d74=1200\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"a\u00d7b=c\" answer cell with the new values,\n# per the commit's diff. Every old string in the mapping is unique within\n# the document, so an exact, case-sensitive Find/Replace is unambiguous.\n$pairs = @(\n  @('2025-06-09 Monday', '2025-06-10 Tuesday'),\n  @('814\u00d79=7326', '880\u00d73=2640'),\n  @('259\u00d76=1554', '831\u00d78=6648'),\n  @('407\u00d78=3256', '687\u00d72=1374'),\n  @('541\u00d75=2705', '621\u00d76=3726'),\n  @('974\u00d79=8766', '474\u00d79=4266'),\n  @('390\u00d77=2730', '834\u00d74=3336'),\n  @('510\u00d76=3060', '707\u00d72=1414'),\n  @('736\u00d72=1472', '414\u00d78=3312'),\n  @('441\u00d77=3087', '861\u00d76=5166'),\n  @('205\u00d78=1640', '950\u00d75=4750'),\n  @('919\u00d73=2757', '408\u00d75=2040'),\n  @('843\u00d76=5058', '628\u00d73=1884'),\n  @('334\u00d79=3006', '801\u00d72=1602'),\n  @('946\u00d77=6622', '755\u00d77=5285'),\n  @('763\u00d76=4578', '575\u00d73=1725'),\n  @('864\u00d75=4320', '797\u00d75=3985'),\n  @('674\u00d75=3370', '533\u00d76=3198'),\n  @('388\u00d73=1164', '750\u00d78=6000'),\n  @('326\u00d72=652', '394\u00d73=1182'),\n  @('498\u00d75=2490', '851\u00d78=6808'),\n  @('125\u00d72=250', '317\u00d74=1268'),\n  @('201\u00d72=402', '983\u00d79=8847'),\n  @('231\u00d76=1386', '218\u00d76=1308'),\n  @('252\u00d77=1764', '261\u00d73=783'),\n  @('856\u00d75=4280', '300\u00d74=1200')\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
